$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.29601514339447
$ws.Range("B1").Value = 2.319894790649414
$ws.Range("C1").Value = 3.148501634597778
$ws.Range("D1").Value = 3.681742668151855
$ws.Range("E1").Value = 1.813809156417847
